$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

function Set-NumCell($cell, $value) {
    $ws.Range($cell).Value = $value
}

# Row 10 - convert A10/D10 from text to numeric, keep other values
Set-NumCell "A10" 29
$ws.Range("B10").Value = "Showers in the Vicinity"
Set-TextCell "C10" "08/05/2024"
Set-NumCell "D10" 13

# Row 11
Set-NumCell "A11" 29
$ws.Range("B11").Value = "Cloudy"
Set-TextCell "C11" "08/05/2024"
Set-NumCell "D11" 14

# Row 12
Set-NumCell "A12" 29
$ws.Range("B12").Value = "Cloudy"
Set-TextCell "C12" "08/05/2024"
Set-NumCell "D12" 15

# Row 13
Set-NumCell "A13" 29
$ws.Range("B13").Value = "Cloudy"
Set-TextCell "C13" "08/05/2024"
Set-NumCell "D13" 16

# Row 14
Set-NumCell "A14" 28
$ws.Range("B14").Value = "Cloudy"
Set-TextCell "C14" "08/05/2024"
Set-NumCell "D14" 17

# Row 15
Set-NumCell "A15" 28
$ws.Range("B15").Value = "Cloudy"
Set-TextCell "C15" "08/06/2024"
Set-NumCell "D15" 8

# Row 16
Set-NumCell "A16" 29
$ws.Range("B16").Value = "Cloudy"
Set-TextCell "C16" "08/06/2024"
Set-NumCell "D16" 9

# Row 17
Set-NumCell "A17" 29
$ws.Range("B17").Value = "Cloudy"
Set-TextCell "C17" "08/06/2024"
Set-NumCell "D17" 10

# Row 18
Set-NumCell "A18" 30
$ws.Range("B18").Value = "Mostly Cloudy"
Set-TextCell "C18" "08/06/2024"
Set-NumCell "D18" 11

# Row 19
Set-NumCell "A19" 31
$ws.Range("B19").Value = "Mostly Cloudy"
Set-TextCell "C19" "08/06/2024"
Set-NumCell "D19" 12

# Row 20 - new row, A20/D20 remain text (like original row 10 was)
Set-TextCell "A20" "29"
$ws.Range("B20").Value = "Rain Shower"
Set-TextCell "C20" "08/06/2024"
Set-TextCell "D20" "13"
